# Update "想去人数" (number of people interested) values for a few events.
# Sheet "展览" (sheet1) and "全部类型" (sheet4) both contain the same rows
# for these events, so both need to be updated to keep the workbook consistent.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 1149
    $ws.Range("F10").Value = 5205
    $ws.Range("F11").Value = 4781
}
